$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9704472422599792
$ws.Range("B1").Value = 1.727505564689636
$ws.Range("C1").Value = 5.780312061309814
$ws.Range("D1").Value = 3.525047302246094
$ws.Range("E1").Value = 0.5851475596427917
